# Final Project Rubric: mark the "ONE PER TEAM" checklist items on the
# "Group rubrik" sheet as completed ("Yes") instead of the placeholder
# prompt text ("completed? Yes or not").
#
# Rows 8-13, column E hold that placeholder for each checklist item
# (Installation instructions, User Manual, Pre-project Team Review,
# Updated Project Proposal, Updated Project Outline, Final architecture
# documentation). Replace all six with "Yes".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group rubrik")
$ws.Activate()

$ws.Range("E8").Value = "Yes"
$ws.Range("E9").Value = "Yes"
$ws.Range("E10").Value = "Yes"
$ws.Range("E11").Value = "Yes"
$ws.Range("E12").Value = "Yes"
$ws.Range("E13").Value = "Yes"

# Match the author's final view state (zoomed in a bit more, cursor left on F12).
$null = $ws.Range("F12").Select()
$excel.ActiveWindow.Zoom = 84
